# Auto-generated edit script: applies per-cell market-data value updates
# across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching
# the upstream "scheduled runner" commit that refreshed computed columns
# H:N (currentAveragePrice*, LevePrice*, LeveProfit*) in the Leve tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64 (sheet ALC), anchor G64=5506
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents() | Out-Null
$ws.Range("N64").ClearContents() | Out-Null

# Row 67 (sheet ALC), anchor G67=5506
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents() | Out-Null
$ws.Range("N67").ClearContents() | Out-Null

# Row 137 (sheet ALC), anchor G137=44013
$ws.Range("H137").Value = 1371.1765
$ws.Range("I137").Value = 1080.4
$ws.Range("K137").Value = 3241.2
$ws.Range("M137").Value = -691.2000000000003

$ws = $wb.Worksheets.Item("ARM")
# Row 30 (sheet ARM), anchor G30=2712
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents() | Out-Null

# Row 32 (sheet ARM), anchor G32=44147
$ws.Range("H32").Value = 4308.952
$ws.Range("I32").Value = 3779.1526
$ws.Range("J32").Value = 12123.5
$ws.Range("K32").Value = 3779.1526
$ws.Range("L32").Value = 12123.5
$ws.Range("M32").Value = -3492.1526
$ws.Range("N32").Value = -12697.5

# Row 115 (sheet ARM), anchor G115=27104
$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -38134

# Row 132 (sheet ARM), anchor G132=43997
$ws.Range("H132").Value = 1540.6333
$ws.Range("I132").Value = 1137.36
$ws.Range("K132").Value = 3412.08
$ws.Range("M132").Value = -882.0799999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (sheet BSM), anchor G86=12526
$ws.Range("H86").Value = 157859.92
$ws.Range("I86").Value = 4568
$ws.Range("J86").Value = 668833
$ws.Range("K86").Value = 4568
$ws.Range("L86").Value = 668833
$ws.Range("M86").Value = -3445
$ws.Range("N86").Value = -671079

# Row 89 (sheet BSM), anchor G89=12526
$ws.Range("H89").Value = 157859.92
$ws.Range("I89").Value = 4568
$ws.Range("J89").Value = 668833
$ws.Range("K89").Value = 22840
$ws.Range("L89").Value = 3344165
$ws.Range("M89").Value = -17224
$ws.Range("N89").Value = -3355397

# Row 134 (sheet BSM), anchor G134=43998
$ws.Range("H134").Value = 4537.4062
$ws.Range("I134").Value = 4920.2593
$ws.Range("K134").Value = 14760.7779
$ws.Range("M134").Value = -12225.7779

# Row 141 (sheet BSM), anchor G141=43278
$ws.Range("H141").Value = 51296.5
$ws.Range("J141").Value = 65926.336
$ws.Range("L141").Value = 65926.336
$ws.Range("N141").Value = -76286.336

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (sheet CRP), anchor G31=44023
$ws.Range("H31").Value = 2332.9092
$ws.Range("I31").Value = 1956.7142
$ws.Range("K31").Value = 1956.7142
$ws.Range("M31").Value = -1661.7142

# Row 34 (sheet CRP), anchor G34=44023
$ws.Range("H34").Value = 2332.9092
$ws.Range("I34").Value = 1956.7142
$ws.Range("K34").Value = 1956.7142
$ws.Range("M34").Value = -1754.7142

# Row 105 (sheet CRP), anchor G105=19928
$ws.Range("H105").Value = 806.4286
$ws.Range("I105").Value = 778.5454999999999
$ws.Range("K105").Value = 778.5454999999999
$ws.Range("M105").Value = 968.4545000000001

# Row 134 (sheet CRP), anchor G134=44020
$ws.Range("H134").Value = 1840.9429
$ws.Range("I134").Value = 1644.8125
$ws.Range("J134").Value = 3933
$ws.Range("K134").Value = 4934.4375
$ws.Range("L134").Value = 11799
$ws.Range("M134").Value = -2399.4375
$ws.Range("N134").Value = -16869

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (sheet CUL), anchor G33=4867
$ws.Range("H33").Value = 133
$ws.Range("I33").Value = 123.333336
$ws.Range("J33").Value = 162
$ws.Range("K33").Value = 740.000016
$ws.Range("L33").Value = 972
$ws.Range("M33").Value = -457.000016
$ws.Range("N33").Value = -1538

# Row 50 (sheet CUL), anchor G50=4725
$ws.Range("H50").Value = 200140670
$ws.Range("I50").Value = 348884.5
$ws.Range("J50").Value = 333335200
$ws.Range("K50").Value = 1046653.5
$ws.Range("L50").Value = 1000005600
$ws.Range("M50").Value = -1046172.5
$ws.Range("N50").Value = -1000006562

# Row 53 (sheet CUL), anchor G53=4725
$ws.Range("H53").Value = 200140670
$ws.Range("I53").Value = 348884.5
$ws.Range("J53").Value = 333335200
$ws.Range("K53").Value = 1046653.5
$ws.Range("L53").Value = 1000005600
$ws.Range("M53").Value = -1046172.5
$ws.Range("N53").Value = -1000006562

# Row 68 (sheet CUL), anchor G68=12895
$ws.Range("H68").Value = 743.4
$ws.Range("J68").Value = 829.25
$ws.Range("L68").Value = 2487.75
$ws.Range("N68").Value = -4109.75

# Row 71 (sheet CUL), anchor G71=12895
$ws.Range("H71").Value = 743.4
$ws.Range("J71").Value = 829.25
$ws.Range("L71").Value = 7463.25
$ws.Range("N71").Value = -15575.25

# Row 88 (sheet CUL), anchor G88=12851
$ws.Range("H88").Value = 4747.875
$ws.Range("I88").Value = 2507
$ws.Range("J88").Value = 5494.8335
$ws.Range("K88").Value = 7521
$ws.Range("L88").Value = 16484.5005
$ws.Range("M88").Value = -7093
$ws.Range("N88").Value = -17340.5005

# Row 91 (sheet CUL), anchor G91=12851
$ws.Range("H91").Value = 4747.875
$ws.Range("I91").Value = 2507
$ws.Range("J91").Value = 5494.8335
$ws.Range("K91").Value = 7521
$ws.Range("L91").Value = 16484.5005
$ws.Range("M91").Value = -6039
$ws.Range("N91").Value = -19448.5005

$ws = $wb.Worksheets.Item("GSM")
# Row 28 (sheet GSM), anchor G28=2063
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents() | Out-Null

# Row 102 (sheet GSM), anchor G102=36169
$ws.Range("H102").Value = 2455.6191
$ws.Range("I102").Value = 2039
$ws.Range("K102").Value = 2039
$ws.Range("M102").Value = -417

# Row 126 (sheet GSM), anchor G126=36184
$ws.Range("H126").Value = 2462022.2
$ws.Range("I126").Value = 3971475.8
$ws.Range("K126").Value = 11914427.4
$ws.Range("M126").Value = -11911957.4

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (sheet LTW), anchor G40=36248
$ws.Range("H40").Value = 12799.947
$ws.Range("I40").Value = 12150.417
$ws.Range("J40").Value = 13913.429
$ws.Range("K40").Value = 12150.417
$ws.Range("L40").Value = 13913.429
$ws.Range("M40").Value = -12014.417
$ws.Range("N40").Value = -14185.429

# Row 82 (sheet LTW), anchor G82=12565
$ws.Range("H82").Value = 2018.9231
$ws.Range("I82").Value = 1925
$ws.Range("J82").Value = 2060.6667
$ws.Range("K82").Value = 1925
$ws.Range("L82").Value = 2060.6667
$ws.Range("M82").Value = -1564
$ws.Range("N82").Value = -2782.6667

# Row 85 (sheet LTW), anchor G85=12565
$ws.Range("H85").Value = 2018.9231
$ws.Range("I85").Value = 1925
$ws.Range("J85").Value = 2060.6667
$ws.Range("K85").Value = 1925
$ws.Range("L85").Value = 2060.6667
$ws.Range("M85").Value = -677
$ws.Range("N85").Value = -4556.6667

# Row 132 (sheet LTW), anchor G132=44058
$ws.Range("H132").Value = 2427.7568
$ws.Range("I132").Value = 1531.2142
$ws.Range("J132").Value = 2973.4783
$ws.Range("K132").Value = 4593.642599999999
$ws.Range("L132").Value = 8920.4349
$ws.Range("M132").Value = -2063.642599999999
$ws.Range("N132").Value = -13980.4349

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (sheet WVR), anchor G122=36208
$ws.Range("H122").Value = 31654.154
$ws.Range("I122").Value = 61352.69
$ws.Range("K122").Value = 184058.07
$ws.Range("M122").Value = -181608.07

# Row 132 (sheet WVR), anchor G132=44029
$ws.Range("H132").Value = 1795.4138
$ws.Range("J132").Value = 3147.5715
$ws.Range("L132").Value = 9442.7145
$ws.Range("N132").Value = -14502.7145

# Row 141 (sheet WVR), anchor G141=42505
$ws.Range("H141").Value = 81874.25
$ws.Range("J141").Value = 81874.25
$ws.Range("L141").Value = 81874.25
$ws.Range("N141").Value = -92234.25
